# Apply updated market-price / profit figures to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 361.84616
$ws.Range("I28").Value = 361.84616
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 361.84616
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 123.15384
$ws.Range("N28").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 732884.3
$ws.Range("I92").Value = 809975.3
$ws.Range("J92").Value = 520
$ws.Range("K92").Value = 809975.3
$ws.Range("L92").Value = 520
$ws.Range("M92").Value = -808727.3
$ws.Range("N92").Value = -3016

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 257.16666
$ws.Range("I101").Value = 192
$ws.Range("J101").Value = 387.5
$ws.Range("K101").Value = 576
$ws.Range("L101").Value = 1162.5
$ws.Range("M101").Value = 1046
$ws.Range("N101").Value = -4406.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 562.925
$ws.Range("I107").Value = 560.2
$ws.Range("J107").Value = 571.1
$ws.Range("K107").Value = 560.2
$ws.Range("L107").Value = 571.1
$ws.Range("M107").Value = 1359.8
$ws.Range("N107").Value = -4411.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 60515.055
$ws.Range("I125").Value = 168900
$ws.Range("J125").Value = 6322.5835
$ws.Range("K125").Value = 1520100
$ws.Range("L125").Value = 56903.2515
$ws.Range("M125").Value = -1517640
$ws.Range("N125").Value = -61823.2515

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 90948160
$ws.Range("J134").Value = 90948160
$ws.Range("L134").Value = 90948160
$ws.Range("N134").Value = -90958300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 65000
$ws.Range("J136").Value = 65000
$ws.Range("L136").Value = 65000
$ws.Range("N136").Value = -75200

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 45988.57
$ws.Range("J139").Value = 45988.57
$ws.Range("L139").Value = 45988.57
$ws.Range("N139").Value = -56268.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3797244.8
$ws.Range("I32").Value = 5279.7734
$ws.Range("K32").Value = 5279.7734
$ws.Range("M32").Value = -4992.7734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 37502396
$ws.Range("I74").Value = 51726684
$ws.Range("J74").Value = 1999.2727
$ws.Range("K74").Value = 51726684
$ws.Range("L74").Value = 1999.2727
$ws.Range("M74").Value = -51725810
$ws.Range("N74").Value = -3747.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 37502396
$ws.Range("I77").Value = 51726684
$ws.Range("J77").Value = 1999.2727
$ws.Range("K77").Value = 258633420
$ws.Range("L77").Value = 9996.363499999999
$ws.Range("M77").Value = -258629052
$ws.Range("N77").Value = -18732.3635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5307.5386
$ws.Range("I107").Value = 5583.3335
$ws.Range("J107").Value = 1998
$ws.Range("K107").Value = 5583.3335
$ws.Range("L107").Value = 1998
$ws.Range("M107").Value = -3663.3335
$ws.Range("N107").Value = -5838

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3713.3677
$ws.Range("I134").Value = 1489.159
$ws.Range("J134").Value = 7791.0835
$ws.Range("K134").Value = 4467.477000000001
$ws.Range("L134").Value = 23373.2505
$ws.Range("M134").Value = -1932.477000000001
$ws.Range("N134").Value = -28443.2505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2972.0625
$ws.Range("J16").Value = 2259.375
$ws.Range("L16").Value = 2259.375
$ws.Range("N16").Value = -2833.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2972.0625
$ws.Range("J113").Value = 2259.375
$ws.Range("L113").Value = 2259.375
$ws.Range("N113").Value = -6599.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 40002030
$ws.Range("I122").Value = 45455400
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 136366200
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -136363750
$ws.Range("N122").Value = -36900.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1815.2593
$ws.Range("I132").Value = 1089.1724
$ws.Range("J132").Value = 2657.52
$ws.Range("K132").Value = 3267.5172
$ws.Range("L132").Value = 7972.559999999999
$ws.Range("M132").Value = -737.5171999999998
$ws.Range("N132").Value = -13032.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 891.55
$ws.Range("I68").Value = 735.7414
$ws.Range("J68").Value = 1106.7142
$ws.Range("K68").Value = 2207.2242
$ws.Range("L68").Value = 3320.1426
$ws.Range("M68").Value = -1396.2242
$ws.Range("N68").Value = -4942.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 891.55
$ws.Range("I71").Value = 735.7414
$ws.Range("J71").Value = 1106.7142
$ws.Range("K71").Value = 6621.6726
$ws.Range("L71").Value = 9960.427799999999
$ws.Range("M71").Value = -2565.6726
$ws.Range("N71").Value = -18072.4278

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3955.2
$ws.Range("I118").Value = 885.6
$ws.Range("J118").Value = 5490
$ws.Range("K118").Value = 2656.8
$ws.Range("L118").Value = 16470
$ws.Range("M118").Value = -1413.8
$ws.Range("N118").Value = -18956

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3168.6667
$ws.Range("I122").Value = 2562.5
$ws.Range("J122").Value = 5593.3335
$ws.Range("K122").Value = 7687.5
$ws.Range("L122").Value = 16780.0005
$ws.Range("M122").Value = -5237.5
$ws.Range("N122").Value = -21680.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 66669350
$ws.Range("I40").Value = 90910184
$ws.Range("K40").Value = 90910184
$ws.Range("M40").Value = -90910048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 186435.73
$ws.Range("I93").Value = 228382.14
$ws.Range("J93").Value = 1871.6
$ws.Range("K93").Value = 228382.14
$ws.Range("L93").Value = 1871.6
$ws.Range("M93").Value = -227134.14
$ws.Range("N93").Value = -4367.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 38505676
$ws.Range("I132").Value = 58887930
$ws.Range("J132").Value = 5866.4443
$ws.Range("K132").Value = 176663790
$ws.Range("L132").Value = 17599.3329
$ws.Range("M132").Value = -176661260
$ws.Range("N132").Value = -22659.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10871131
$ws.Range("I136").Value = 16130086
$ws.Range("J136").Value = 2624
$ws.Range("K136").Value = 48390258
$ws.Range("L136").Value = 7872
$ws.Range("M136").Value = -48387708
$ws.Range("N136").Value = -12972

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 562.44446
$ws.Range("I100").Value = 567.8333
$ws.Range("K100").Value = 1135.6666
$ws.Range("M100").Value = -594.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 23590
$ws.Range("J108").Value = 23590
$ws.Range("L108").Value = 23590
$ws.Range("N108").Value = -31270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12559.462
$ws.Range("I122").Value = 15473.111
$ws.Range("J122").Value = 6003.75
$ws.Range("K122").Value = 46419.333
$ws.Range("L122").Value = 18011.25
$ws.Range("M122").Value = -43969.333
$ws.Range("N122").Value = -22911.25
